$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 43 and 44 swap coin identity (name/link) and get refreshed price/volume values.
# D43/D44 are numeric-looking strings ("12.50"/"16.14") -> force Text so they stay
# as literal strings (matching the source data's text-typed Price column) instead of
# being auto-coerced to numbers.
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.50"
$ws.Range("E43").Value = "  +0.42%  "

$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.14"
$ws.Range("E44").Value = "  +4.12%  "

# Price / 1h-volume refresh across the remaining rows.
# Price (column D) values that look like plain decimal numbers need NumberFormat
# forced to Text ("@") before the assignment so Excel doesn't silently convert
# them to numeric cells; values containing two dots (e.g. "23.446.72") are never
# parsed as numbers so no special handling is required there.
$ws.Range("D2").Value = "23.446.72"
$ws.Range("E2").Value = "  +1.13%  "
$ws.Range("D3").Value = "1.638.42"
$ws.Range("E3").Value = "  +2.36%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "306.72"
$ws.Range("E6").Value = "  +1.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3766"
$ws.Range("E7").Value = "  -0.41%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "52.27"
$ws.Range("E8").Value = "  +0.23%  "
$ws.Range("E9").Value = "  +0.81%  "
$ws.Range("E10").Value = "  -0.09%  "
$ws.Range("E11").Value = "  +0.42%  "
$ws.Range("E12").Value = "  +0.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.89"
$ws.Range("E13").Value = "  +1.05%  "
$ws.Range("E14").Value = "  +0.77%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001274"
$ws.Range("E15").Value = "  +2.33%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.362"
$ws.Range("E16").Value = "  -0.48%  "
$ws.Range("D17").Value = "1.640.84"
$ws.Range("E17").Value = "  +2.61%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.57"
$ws.Range("E18").Value = "  +0.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06932"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.13"
$ws.Range("E20").Value = "  +0.35%  "
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").Value = "23.445.13"
$ws.Range("E23").Value = "  +1.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.76"
$ws.Range("E24").Value = "  -1.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.100"
$ws.Range("E25").Value = "  +3.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.423"
$ws.Range("E26").Value = "  +1.78%  "
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "150.68"
$ws.Range("E28").Value = "  +0.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.367"
$ws.Range("E29").Value = "  +2.14%  "
$ws.Range("E30").Value = "  +0.98%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.309"
$ws.Range("E31").Value = "  -3.22%  "
$ws.Range("D32").Value = "1.823.91"
$ws.Range("E32").Value = "  +2.47%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.787"
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9659"
$ws.Range("E34").Value = "  -0.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02817"
$ws.Range("E35").Value = "  +4.08%  "
$ws.Range("E36").Value = "  +0.47%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.07300"
$ws.Range("E37").Value = "  -2.58%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2522"
$ws.Range("E38").Value = "  +0.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.08838"
$ws.Range("E39").Value = "  +0.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.114"
$ws.Range("E40").Value = "  +0.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.375"
$ws.Range("E41").Value = "  +0.97%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7090"
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6534"
$ws.Range("E45").Value = "  +0.16%  "
$ws.Range("E46").Value = "  +1.08%  "
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("E48").Value = "  +0.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07958"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "128.82"
$ws.Range("E50").Value = "  -2.44%  "
$ws.Range("E51").Value = "  +0.25%  "
